$p = $ppt.ActivePresentation

# --- Slide 3: update the title text (merge the two runs into one) ---
$s3 = $p.Slides.Item(3)
$title = $s3.Shapes.Item(1)
# First collapse to a placeholder string so the COM text-diffing engine
# doesn't try to preserve the old runs' split point, then set the final
# text -- this yields a single <a:r> run like a full retype would.
$title.TextFrame.TextRange.Text = "_"
$title.TextFrame.TextRange.Text = "Household Food Security by Income, 2021"

# --- Slide 5: change the table's style id ---
$s5 = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(3).Table
$tbl.ApplyStyle("{0DB4F5D1-2DAD-44FD-89DC-8E78F9413ADD}")
